# Fill in the presentation date on the cover page.
#
# Before:  "Présenté le                                       par :"   (one run)
# After:   "Présenté le " + "15 mars 2024" + "  par :"                 (three runs)
#
# The middle run carries the newly-typed date; the outer two runs are what is
# left of the original run once the date was typed in the middle of it. All
# three runs keep the original run's formatting (Latin Modern Math, bold, 28pt).

$d = $word.ActiveDocument

# Locate the cover-page paragraph holding the "Présenté le ... par :" line.
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "Présenté le*par :*") {
        $target = $p.Range
        break
    }
}

if ($target -eq $null) {
    throw "Could not find the 'Présenté le ... par :' paragraph"
}

# Range covering just the paragraph's text, without its trailing paragraph mark.
$r = $d.Range($target.Start, $target.End - 1)

# Rewrite that range as three runs (same rPr, split text) via InsertXML so the
# saved document ends up with distinct <w:r> elements instead of one merged run.
$rPr = '<w:rPr><w:rFonts w:ascii="Latin Modern Math" w:hAnsi="Latin Modern Math"/><w:b/><w:bCs/><w:sz w:val="28"/><w:szCs w:val="28"/></w:rPr>'

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
         '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
           '<pkg:xmlData>' +
             '<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
               '<w:body>' +
                 '<w:p>' +
                   '<w:r>' + $rPr + '<w:t xml:space="preserve">Présenté le </w:t></w:r>' +
                   '<w:r>' + $rPr + '<w:t>15 mars 2024</w:t></w:r>' +
                   '<w:r>' + $rPr + '<w:t xml:space="preserve">  par :</w:t></w:r>' +
                 '</w:p>' +
               '</w:body>' +
             '</w:document>' +
           '</pkg:xmlData>' +
         '</pkg:part>' +
       '</pkg:package>'

$r.InsertXML($xml)
